# Trade #23 closed at 2026-02-17 12:37:20 - unknown UNKNOWN +0.000%
#
# Adds the newly-closed trade (#23) to the "All Trades" and "MarketMaking"
# log sheets, and rolls the aggregate statistics forward on the "Summary"
# and "Strategy Status" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - roll aggregate stats forward to include trade #23
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.64   # Current Capital
$summary.Range("B4").Value = 0.64      # Total P&L $
$summary.Range("B5").Value = 0.54 + 0.02  # Total P&L %  -> 0.5600000000000001
$summary.Range("B6").Value = 23        # Total Trades
$summary.Range("B7").Value = 9         # Winning Trades
$summary.Range("B9").Value = 39.13     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.64
$status.Range("D4").Value = 23
$status.Range("E4").Value = 0.64
$status.Range("F4").Value = 0.64
$status.Range("G4").Value = 39.13

# ---------------------------------------------------------------------
# Append the new trade row (row 24) to both trade-log sheets
# ---------------------------------------------------------------------
$tradeSheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 24

    $ws.Cells.Item($row, 1).Value = 23

    # Keep the date as plain text (matches the rest of the column) instead
    # of letting automatic date recognition turn it into a date serial.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"

    $ws.Cells.Item($row, 3).Value = "12:37:13"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.57
    $ws.Cells.Item($row, 7).Value = 0.62
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 8.7719
    $ws.Cells.Item($row, 10).Value = 0.05
    $ws.Cells.Item($row, 11).Value = 100.64
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}
